$d = $word.ActiveDocument
$W = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Replace the whole text of the (single, first-encountered) paragraph that
# contains $find with the literal paragraph XML given in $xml. Using
# Range.InsertXML (rather than Find-replace or Range.Text=) preserves the
# paragraph's other runs (e.g. a leading empty <w:r/>) exactly as-is.
function Set-ParagraphXml($find, $xml) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($find)
    if (-not $ok) {
        Write-Host "NOT FOUND: $find"
        return
    }
    $para = $rng.Paragraphs(1)
    $para.Range.InsertXML($xml)
}

# Heading1 title at the top of the document (no leading empty run).
Set-ParagraphXml `
    "Play Fetching Fruits for Free - A Classic 5x3 Fruit-Themed Slot" `
    "<w:p $W><w:pPr><w:pStyle w:val=""Heading1""/></w:pPr><w:r><w:t>Play Fetching Fruits Free - Classic Fruit-Themed Slot Game</w:t></w:r></w:p>"

# "What we like" bullet points (reordered + reworded).
Set-ParagraphXml `
    "Pleasant musical component with electronic sounds and flutes" `
    "<w:p $W><w:pPr><w:pStyle w:val=""ListBullet""/><w:spacing w:line=""240"" w:lineRule=""auto""/><w:ind w:left=""720""/></w:pPr><w:r/><w:r><w:t>Classic fruit-themed slot game</w:t></w:r></w:p>"

Set-ParagraphXml `
    "Symbols and paylines are clearly visible and straightforward" `
    "<w:p $W><w:pPr><w:pStyle w:val=""ListBullet""/><w:spacing w:line=""240"" w:lineRule=""auto""/><w:ind w:left=""720""/></w:pPr><w:r/><w:r><w:t>Pleasant musical component</w:t></w:r></w:p>"

Set-ParagraphXml `
    "Autoplay feature allows for a more hands-off approach to playing slots" `
    "<w:p $W><w:pPr><w:pStyle w:val=""ListBullet""/><w:spacing w:line=""240"" w:lineRule=""auto""/><w:ind w:left=""720""/></w:pPr><w:r/><w:r><w:t>Autoplay feature for convenient gameplay</w:t></w:r></w:p>"

Set-ParagraphXml `
    "Higher bet amounts increase chances of winning and potential winnings" `
    "<w:p $W><w:pPr><w:pStyle w:val=""ListBullet""/><w:spacing w:line=""240"" w:lineRule=""auto""/><w:ind w:left=""720""/></w:pPr><w:r/><w:r><w:t>Clear and straightforward symbols and paylines</w:t></w:r></w:p>"

# "What we don't like" bullet points (reworded).
Set-ParagraphXml `
    "Bare internal menu lacks information on volatility and RTP" `
    "<w:p $W><w:pPr><w:pStyle w:val=""ListBullet""/><w:spacing w:line=""240"" w:lineRule=""auto""/><w:ind w:left=""720""/></w:pPr><w:r/><w:r><w:t>Limited information in the internal menu</w:t></w:r></w:p>"

Set-ParagraphXml `
    "Limited number of paylines (only 5)" `
    "<w:p $W><w:pPr><w:pStyle w:val=""ListBullet""/><w:spacing w:line=""240"" w:lineRule=""auto""/><w:ind w:left=""720""/></w:pPr><w:r/><w:r><w:t>Lack of volatility and RTP details</w:t></w:r></w:p>"

# Bold "title" run and italic "meta description" run near the very end of
# the document (each paragraph has a leading empty run too).
Set-ParagraphXml `
    "Play Fetching Fruits for Free - A Classic 5x3 Fruit-Themed Slot" `
    "<w:p $W><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fetching Fruits Free - Classic Fruit-Themed Slot Game</w:t></w:r></w:p>"

Set-ParagraphXml `
    "Explore this classic 5x3 fruit-themed slot, Fetching Fruits by NewAge Games. Enjoy pleasant music and high-level symbols with win potential. Play for free now!" `
    "<w:p $W><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Fetching Fruits and play this classic fruit-themed slot game for free.</w:t></w:r></w:p>"
